# MigrationRenamer.xlsx update:
#  - Previous run's target filenames (2020_11_06_...) become this run's source "Files"
#    (done by replacing the 2020_10_05_ prefix with 2020_11_06_ in the existing rows)
#  - Two new migrations show up: "create_masters_table" (inserted in the middle, after
#    create_media_table) and "create_remote_items_table" (appended at the end)
#  - The calculated "Date" column now stamps rows with 2021_02_08_ instead of 2020_11_06_

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Step 1: the 19 existing data rows already hold last run's *destination* names;
# turn them into this run's *source* names by bumping the embedded date prefix.
$ws.Range("A2:A20").Replace("2020_10_05_", "2020_11_06_")

# Step 2: make room for the new "create_masters_table" migration right after
# create_media_table (i.e. before the old row 5 / create_user_logins_table).
$ws.Rows.Item(5).Insert()

# Step 3: grow the table to cover the two extra rows (21 data rows total + header).
$lo.Resize($ws.Range("A1:F22"))

# Step 4: fill in the "Files" (source) column for the two brand-new rows.
$ws.Range("A5").Value = "2020_11_06_000020_create_masters_table.php"
$ws.Range("A22").Value = "2020_11_06_000021_create_remote_items_table.php"

# Step 5: (re)build all of the calculated columns for every data row so that the
# whole table is internally consistent, exactly mirroring the table's own
# calculated-column formulas.
$ws.Range("B2:B22").Formula = "=MID(Table1[[#This Row],[Files]],18,LEN(Table1[[#This Row],[Files]]))"
$ws.Range("C2:C22").Formula = "=`"2021_02_08_`""
$ws.Range("D2:D22").Formula = "=REPT(`"0`",6-LEN(MATCH(Table1[[#This Row],[Name]],Table1[Name],0)))&MATCH(Table1[[#This Row],[Name]],Table1[Name],0)"
$ws.Range("E2:E22").Formula = "=Table1[[#This Row],[Date]]&Table1[[#This Row],[Seq]]&Table1[[#This Row],[Name]]"
$ws.Range("F2:F22").Formula = "=`"ren `"&Table1[[#This Row],[Files]]&`" `"&Table1[[#This Row],[Replace]]&`";`""

# Step 6: refresh the dimension / selection the way Excel would after this edit.
$ws.Range("F2:F22").Select()

$wb.Save()
